# Add the second plate's multiplex data file (1/14/20) to the map:
# every data row's multiplex_file column now points at the new plate's
# workbook instead of the old placeholder name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2:E11").Value = "C:\Users\Ryan\OneDrive - Iota Bio\data\multiplex\1_14_20_rat_plate1.xlsx"

$ws.Range("D16").Select()

$wb.Save()
